$wb = $excel.ActiveWorkbook

# --- "Produtos" sheet: replace the "laptops" sample row with the ---
# --- "LAPTOPS" / "HP PAVILION 15T TOUCH LAPTOP" product, and add  ---
# --- a new "Geral" / generic-scroll validation row.               ---
$wsProdutos = $wb.Worksheets.Item("Produtos")

$wsProdutos.Range("A5").Value = "LAPTOPS"
$wsProdutos.Range("B5").Value = "HP PAVILION 15T TOUCH LAPTOP"
$wsProdutos.Range("C5").Value = "HP PAVILION 15T TOUCH LAPTOP"

$wsProdutos.Range("A15").Value = "Geral"
$wsProdutos.Range("B15").Value = "No results for"
$wsProdutos.Range("C15").Value = "C(14,0) v(14,1)"

$wsProdutos.Range("A12").Select() | Out-Null

# --- "Cadastro" sheet: bump the test-run/build tag, and restore it ---
# --- as the active tab (it was active before this edit too).       ---
$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsCadastro.Range("B2").Value = "Wilkerbn503"
$wsCadastro.Select() | Out-Null
